$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the CasesTab query (B2): drop the cohort-related MATCH/column now that
# the extra browser-opening cohort lookup was removed from the script.
$newB2 = @"
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)

MATCH (c)<--(diag:diagnosis)
WHERE s.clinical_study_designation IN ['UBC01'] and demo.breed in ['Mixed Breed', 'Scottish Terrier','Shetland Sheepdog']and diag.disease_term in ['Bladder Cancer','Healthy Control'] and diag.primary_disease_site in ['Bladder', 'Bladder, Prostate', 'Bladder, Urethra', 'Bladder, Urethra, Prostate', 'Not Applicable']
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS ``Case ID`` ,
        coalesce(s.clinical_study_designation, '') AS ``Study Code`` ,
        coalesce(s.clinical_study_type, '') AS  ``Study Type``,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS ``Stage of Disease`` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS ``Neutered Status``,
        coalesce(demo.weight, '') AS ``Weight (kg)``,
        coalesce(diag.best_response, '') AS ``Response to Treatment``
"@

$ws.Range("B2").Value = $newB2

# Row heights settle slightly differently once the shorter query text is in
# place (and FilesTab/SamplesTab rows re-wrap too).
$ws.Rows(2).RowHeight = 315
$ws.Rows(3).RowHeight = 315
$ws.Rows(4).RowHeight = 300

# Selection ends up on B2 (the cell that was edited) instead of C4.
$ws.Range("B2").Select()

Write-Output "edit complete"
